# Auto-generated edit script applying scheduled market-data refresh
# to the Ultros_Profits workbook (columns H-N per sheet/row).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 83.375
$ws.Range("I8").Value = 23.857143
$ws.Range("K8").Value = 71.57142899999999
$ws.Range("M8").Value = 67.42857100000001
# Row 9
$ws.Range("H9").Value = 217.625
$ws.Range("I9").Value = 209.8
$ws.Range("J9").Value = 221.18182
$ws.Range("K9").Value = 209.8
$ws.Range("L9").Value = 221.18182
$ws.Range("M9").Value = -40.80000000000001
$ws.Range("N9").Value = -559.18182
# Row 40
$ws.Range("H40").Value = 3913.5454
$ws.Range("I40").Value = 2592.8572
$ws.Range("K40").Value = 2592.8572
$ws.Range("M40").Value = -2417.8572
# Row 42
$ws.Range("H42").Value = 1799
$ws.Range("I42").Value = 197.57143
$ws.Range("J42").Value = 13009
$ws.Range("K42").Value = 592.71429
$ws.Range("L42").Value = 39027
$ws.Range("M42").Value = -362.71429
$ws.Range("N42").Value = -39487
# Row 61
$ws.Range("H61").Value = 257.66666
$ws.Range("I61").Value = 242.375
$ws.Range("K61").Value = 727.125
$ws.Range("M61").Value = -555.125
# Row 80
$ws.Range("H80").Value = 2649.8125
$ws.Range("I80").Value = 997.25
$ws.Range("K80").Value = 2991.75
$ws.Range("M80").Value = -1993.75
# Row 83
$ws.Range("H83").Value = 2649.8125
$ws.Range("I83").Value = 997.25
$ws.Range("K83").Value = 8975.25
$ws.Range("M83").Value = -3983.25
# Row 112
$ws.Range("H112").Value = 1327.1154
$ws.Range("J112").Value = 1326.3043
$ws.Range("L112").Value = 3978.9129
$ws.Range("N112").Value = -6194.9129
# Row 113
$ws.Range("H113").Value = 16701.666
$ws.Range("I113").Value = 15000
$ws.Range("J113").Value = 16856.363
$ws.Range("K113").Value = 15000
$ws.Range("L113").Value = 16856.363
$ws.Range("M113").Value = -11746
$ws.Range("N113").Value = -23364.363
# Row 135
$ws.Range("H135").Value = 3209.9062
$ws.Range("I135").Value = 3048.2
$ws.Range("J135").Value = 3352.5881
$ws.Range("K135").Value = 27433.8
$ws.Range("L135").Value = 30173.2929
$ws.Range("M135").Value = -24898.8
$ws.Range("N135").Value = -35243.2929
# Row 141
$ws.Range("H141").Value = 5469.96
$ws.Range("I141").Value = 4815.174
$ws.Range("K141").Value = 14445.522
$ws.Range("M141").Value = -9265.522000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 19493.9
$ws.Range("I2").Value = 23826.732
$ws.Range("J2").Value = 6495.4
$ws.Range("K2").Value = 23826.732
$ws.Range("L2").Value = 6495.4
$ws.Range("M2").Value = -23713.732
$ws.Range("N2").Value = -6721.4
# Row 32
$ws.Range("H32").Value = 4616.625
$ws.Range("I32").Value = 4269.918
$ws.Range("J32").Value = 11666.333
$ws.Range("K32").Value = 4269.918
$ws.Range("L32").Value = 11666.333
$ws.Range("M32").Value = -3982.918
$ws.Range("N32").Value = -12240.333
# Row 45
$ws.Range("H45").Value = 8068.231
$ws.Range("I45").Value = 4498
$ws.Range("J45").Value = 9655
$ws.Range("K45").Value = 4498
$ws.Range("L45").Value = 9655
$ws.Range("M45").Value = -4121
$ws.Range("N45").Value = -10409
# Row 61
$ws.Range("H61").Value = 18675.625
$ws.Range("I61").Value = 2531
$ws.Range("J61").Value = 45583.332
$ws.Range("K61").Value = 2531
$ws.Range("L61").Value = 45583.332
$ws.Range("M61").Value = -2319
$ws.Range("N61").Value = -46007.332
# Row 102
$ws.Range("H102").Value = 17550654
$ws.Range("I102").Value = 4723.1875
$ws.Range("K102").Value = 4723.1875
$ws.Range("M102").Value = -3101.1875
# Row 116
$ws.Range("H116").Value = 19493.9
$ws.Range("I116").Value = 23826.732
$ws.Range("J116").Value = 6495.4
$ws.Range("K116").Value = 23826.732
$ws.Range("L116").Value = 6495.4
$ws.Range("M116").Value = -21532.732
$ws.Range("N116").Value = -11083.4
# Row 132
$ws.Range("H132").Value = 4098
$ws.Range("J132").Value = 8666.666999999999
$ws.Range("L132").Value = 26000.001
$ws.Range("N132").Value = -31060.001
# Row 136
$ws.Range("H136").Value = 18675.625
$ws.Range("I136").Value = 2531
$ws.Range("J136").Value = 45583.332
$ws.Range("K136").Value = 7593
$ws.Range("L136").Value = 136749.996
$ws.Range("M136").Value = -5043
$ws.Range("N136").Value = -141849.996

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 19493.9
$ws.Range("I3").Value = 23826.732
$ws.Range("J3").Value = 6495.4
$ws.Range("K3").Value = 23826.732
$ws.Range("L3").Value = 6495.4
$ws.Range("M3").Value = -23712.732
$ws.Range("N3").Value = -6723.4
# Row 20
$ws.Range("H20").Value = 5393.7896
$ws.Range("I20").Value = 4029
$ws.Range("K20").Value = 4029
$ws.Range("M20").Value = -3782
# Row 94
$ws.Range("H94").Value = 2366784.5
$ws.Range("I94").Value = 1947.3055
$ws.Range("J94").Value = 10106252
$ws.Range("K94").Value = 1947.3055
$ws.Range("L94").Value = 10106252
$ws.Range("M94").Value = -1496.3055
$ws.Range("N94").Value = -10107154
# Row 134
$ws.Range("H134").Value = 2675
$ws.Range("I134").Value = 1300
$ws.Range("J134").Value = 4050
$ws.Range("K134").Value = 3900
$ws.Range("L134").Value = 12150
$ws.Range("M134").Value = -1365
$ws.Range("N134").Value = -17220

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 1375
$ws.Range("I25").Value = 666.6667
$ws.Range("K25").Value = 666.6667
$ws.Range("M25").Value = -492.6667
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents() | Out-Null
$ws.Range("N51").ClearContents() | Out-Null
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("N61").ClearContents() | Out-Null
# Row 119
$ws.Range("H119").Value = 31761
$ws.Range("J119").Value = 31761
$ws.Range("L119").Value = 31761
$ws.Range("N119").Value = -41437
# Row 134
$ws.Range("H134").Value = 5059.9653
$ws.Range("I134").Value = 5049.04
$ws.Range("K134").Value = 15147.12
$ws.Range("M134").Value = -12612.12

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 1725.75
$ws.Range("I98").Value = 1600.7142
$ws.Range("J98").Value = 1900.8
$ws.Range("K98").Value = 4802.142599999999
$ws.Range("L98").Value = 5702.4
$ws.Range("M98").Value = -3304.142599999999
$ws.Range("N98").Value = -8698.4
# Row 132
$ws.Range("H132").Value = 1850
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 21615.785
$ws.Range("I43").Value = 8717.286
$ws.Range("K43").Value = 8717.286
$ws.Range("M43").Value = -8566.286
# Row 102
$ws.Range("H102").Value = 4094.8
$ws.Range("I102").Value = 3588
$ws.Range("K102").Value = 3588
$ws.Range("M102").Value = -1966
# Row 113
$ws.Range("H113").Value = 12809.392
$ws.Range("I113").Value = 8549.875
$ws.Range("K113").Value = 8549.875
$ws.Range("M113").Value = -6379.875
# Row 132
$ws.Range("H132").Value = 11940.588
$ws.Range("I132").Value = 10915.833
$ws.Range("J132").Value = 14400
$ws.Range("K132").Value = 32747.499
$ws.Range("L132").Value = 43200
$ws.Range("M132").Value = -30217.499
$ws.Range("N132").Value = -48260

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 46053.5
$ws.Range("I25").Value = 46053.5
$ws.Range("K25").Value = 46053.5
$ws.Range("M25").Value = -45823.5
# Row 68
$ws.Range("H68").Value = 3844.4443
$ws.Range("J68").Value = 3933.3333
$ws.Range("L68").Value = 3933.3333
$ws.Range("N68").Value = -5431.3333
# Row 71
$ws.Range("H71").Value = 3844.4443
$ws.Range("J71").Value = 3933.3333
$ws.Range("L71").Value = 19666.6665
$ws.Range("N71").Value = -27154.6665
# Row 93
$ws.Range("H93").Value = 5559941
$ws.Range("I93").Value = 2250.5652
$ws.Range("K93").Value = 2250.5652
$ws.Range("M93").Value = -1002.5652

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 5000.6665
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 5000.6665
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 5000.6665
$ws.Range("M8").ClearContents() | Out-Null
$ws.Range("N8").Value = -5280.6665
# Row 62
$ws.Range("H62").Value = 20837340
# Row 65
$ws.Range("H65").Value = 20837340
